# CaseStudy11 - Paste into a filtered range using the x-marker + sort trick
$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)   # "Yêu cầu"
$ws2 = $wb.Worksheets.Item(2)   # "Cách làm"

# Step 1-2: mark the rows that were visible under the "band 2" filter with an
# "x" in column B, and paste the looked-up values (336 -> 236, 8) from column F
# into column D for those same rows.
$ws2.Range("B6").Value = "x"
$ws2.Range("D6").Value = 236
$ws2.Range("B10").Value = "x"
$ws2.Range("D10").Value = 8

# Step 3: sort A5:F11 by column A (the original row numbers) so everything
# lands back in its original order once the filter is removed.
$sortObj = $ws2.Sort
$sortObj.SortFields.Clear()
$sortObj.SortFields.Add($ws2.Range("A5:A11"))
$sortObj.SetRange($ws2.Range("A5:F11"))
$sortObj.Header = 0
$sortObj.Apply()

# Leave an AutoFilter on the Band/Giá trị table (C4:D11).
$ws2.Range("C4:D11").AutoFilter()

# Repoint the (hidden) _FilterDatabase defined name at the Band/Giá trị
# columns only (C4:D11) instead of the old A4:D11.
$n = $wb.Names.Item(1)
$n.RefersTo = "='Cách làm'!`$C`$4:`$D`$11"

# Restore the on-screen selections as left by the author.
$ws1.Range("B4").Select()
$ws2.Range("B6:F10").Select()
